# Afficher_l'historique_des_conversations.xlsx edit
# - Update the text of several "user story" cells in column C
# - Increase the row height of the table rows (2-10) to 24pt (custom height)
# - Add a bottom border under row 7 (Scenario row) to close off that block
# - Update the sheet zoom level and the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the wording of the user story -------------------------------
$ws.Range("C3").Value  = "utilisateur connecté"
$ws.Range("C4").Value  = "afficher l'historique des conversations d'un contact"
$ws.Range("C5").Value  = "de voir les messages que j'ai envoyé et reçu de ce contact"
$ws.Range("C7").Value  = "L'utilisateur affiche l'historique des conversations d'un de ses contacts"
$ws.Range("C9").Value  = "je clique sur le un bouton de chat"
$ws.Range("C10").Value = "la page de chat s'affiche et je peux voir les messages déjà envoyés"

# --- 2. Grow the table rows a bit (19/20/21pt -> 24pt, custom height) ------
$ws.Range("B2:C10").RowHeight = 24

# --- 3. Close the "Scenario" row (7) with a bottom border ------------------
$scenarioRow = $ws.Range("B7:C7")
$scenarioRow.Borders.Item(9).Weight = -4138

# --- 4. Zoom in a bit and move the selection to C8 --------------------------
$excel.ActiveWindow.Zoom = 130
$ws.Range("C8").Select()
